$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the three "Finger (16) unit-cap" literal measurements that
#    shifted slightly (comparable-layout re-measurement).
# ---------------------------------------------------------------------
$ws.Range("I9").Value  = 32.51
$ws.Range("I10").Value = 31.19
$ws.Range("I11").Value = 30.59

# ---------------------------------------------------------------------
# 2. Turn the "C TOP-BOT" G/H/I cells (rows 15-18) into live formulas
#    referencing the measured cap values above each respective bin
#    (bin caps are now computed on the bottom instead of being pasted
#    as static numbers).
# ---------------------------------------------------------------------
$ws.Range("G15").Formula = "=1000*G8/48"
$ws.Range("H15").Formula = "=1000*H8/48"
$ws.Range("I15").Formula = "=1000*I8/48"

$ws.Range("G16").Formula = "=1000*G9/52"
$ws.Range("H16").Formula = "=1000*H9/52"
$ws.Range("I16").Formula = "=1000*I9/52"

$ws.Range("G17").Formula = "=1000*G10/50"
$ws.Range("H17").Formula = "=1000*H10/50"
$ws.Range("I17").Formula = "=1000*I10/50"

$ws.Range("G18").Formula = "=1000*G11/49"
$ws.Range("H18").Formula = "=1000*H11/49"
$ws.Range("I18").Formula = "=1000*I11/49"

# ---------------------------------------------------------------------
# 3. Extend the legend row (row 28): add "high var" / "higher area/cap"
#    callouts and change the existing "low cap" callouts to
#    "low cap\nhigh var" (wrapped across two lines), bin caps now
#    explained along the bottom legend row.
# ---------------------------------------------------------------------
$ws.Rows(28).RowHeight = 30

$lowCapHighVar = "low cap" + [char]10 + "high var"
$ws.Range("F28").Value = $lowCapHighVar
$ws.Range("G28").Value = $lowCapHighVar
$wrapped = $ws.Range("F28,G28")
$wrapped.Font.Color = 255
$wrapped.NumberFormat = "0.00""aF"""
$wrapped.HorizontalAlignment = -4152
$wrapped.WrapText = $true

$ws.Range("D28").Value = "high var"
$ws.Range("D28").Font.Color = 255
$ws.Range("D28").NumberFormat = "0.00""aF"""
$ws.Range("D28").HorizontalAlignment = -4152

$ws.Range("H28").Value = "higher area/cap"
$ws.Range("H28").Font.Color = 255
$ws.Range("H28").NumberFormat = "0.00""aF"""
$ws.Range("H28").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 4. Highlight the comparable topology columns (Waffle A (16),
#    Waffle (8) A, Finger (8) F) in the header row with bold + green
#    fill, making the layout comparable at a glance.
# ---------------------------------------------------------------------
$headerHighlight = $ws.Range("B1,C1,I1")
$headerHighlight.Font.Bold = $true
$headerHighlight.Interior.Color = 5296274

# ---------------------------------------------------------------------
# 5. Leave the selection where the author left it while doing this
#    comparison (I28 — the new "higher area/cap" legend cell).
# ---------------------------------------------------------------------
$ws.Range("I28").Select()
